# [TECH] Removed password from the Excelsheet | Purushotham
#
# The "Credentials" sheet stored a real password ("mckinsey@123") in
# plain text. Replace it with a placeholder ("xyz") and leave the
# "Credentials" sheet active/selected (at cell B3), matching where the
# author was last working when the secret was scrubbed.

$wb = $excel.ActiveWorkbook

$wsTimesheet   = $wb.Worksheets.Item("Timesheet")
$wsCredentials = $wb.Worksheets.Item("Credentials")

# Scrub the plaintext password value.
$wsCredentials.Range("B2").Value = "xyz"

# Make the Credentials sheet the active tab, with B3 selected - this
# also clears the previous tabSelected/selection state on Timesheet.
$wsCredentials.Activate()
$wsCredentials.Range("B3").Select()
